$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header in column H, matching the formatting of the
# other header cells in row 1 (bold, centered, bordered style).
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add the corresponding data value for the new column in row 2.
$ws.Range("H2").Value = 0
